# "april 19 matches updated"
# Match 34 (LSG vs CSK) DRS review rows appended to the bottom of the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$newRows = @(
    @(34, "LSG", "CSK", 1, "CSK", "LSG", 19, "CSK", "AK Chaudhary", "AKC", "Wide",   "Not Called", "Called",     "MS Dhoni",  "Mohsin Khan",  "Successful",   "No"),
    @(34, "LSG", "CSK", 1, "CSK", "LSG", 19, "LSG", "AK Chaudhary", "AKC", "Wide",   "Called",     "Not Called", "RA Jadeja", "Mohsin Khan",  "Successful",   "No"),
    @(34, "LSG", "CSK", 2, "LSG", "CSK", 18, "CSK", "R Pandit",     "RP",  "Wicket", "Not Out",    "Not Out",    "N Pooran",  "M Pathirana",  "Unsuccessful", "No"),
    @(34, "LSG", "CSK", 2, "LSG", "CSK", 19, "CSK", "AK Chaudhary", "AKC", "NoBall", "Called",     "Called",     "N Pooran",  "TU Deshpande", "Unsuccessful", "No")
)

$startRow = 131
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRows[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowValues[$j]
    }
}

# Leave the cursor where the author left it after typing the last row.
[void]$ws.Range("L136").Select()
